# Add struct für output & sensors
# Populates the "Todo" sheet with the Sensor / Output / Input column
# headers (HTML name-mapping block) plus a big spacer cell, matching the
# committed worksheet layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Todo")

# --- cell values -----------------------------------------------------
# Order matters: new shared-string entries are appended in first-seen
# order, so writes are sequenced to reproduce the original string table.

$ws.Range("A3").Value = "HTML side for name mapping Sensors"
$ws.Range("E3").Value = "HTML side for name mapping digital inputs"
$ws.Range("C3").Value = "HTML side for name mapping digital outputs"

$ws.Range("A4").Value = "MQTT active"
$ws.Range("C4").Value = "MQTT active"
$ws.Range("E4").Value = "MQTT active"

$ws.Range("C5").Value = "current state of Output (on/off/auto)"
$ws.Range("C6").Value = "Output name"

$ws.Range("E5").Value = "current state of Input (on/off)"
$ws.Range("E6").Value = "Input name"

$ws.Range("C7").Value = "Equal to Input"

$ws.Range("A5").Value = "current state of Sensor (Value in °C)"
$ws.Range("A6").Value = "Sensor name"
$ws.Range("A7").Value = "Sensor address"

# --- column widths (best effort bestFit sizing) -----------------------
$ws.Columns.Item(1).ColumnWidth = 40.16666666666666
$ws.Columns.Item(3).ColumnWidth = 39.73697916666666
$ws.Columns.Item(5).ColumnWidth = 38.45182291666666

# --- spacer row with larger heading font ------------------------------
$ws.Range("A9").Font.Size = 18
$ws.Rows.Item(9).RowHeight = 23.25

# --- sheet view / selection -------------------------------------------
[void]$ws.Range("A8").Select()

# --- page setup ---------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
